# Update the "想去人数" (want-to-go count) figures on the 展览 and 全部类型
# sheets to reflect the regenerated data snapshot:
#   F2: 1565 -> 1566
#   F3:   91 ->   92
#   F4:   27 ->   28

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1566
    $ws.Range("F3").Value = 92
    $ws.Range("F4").Value = 28
}
